# Add "Episode 11" results sheet to AYTO VIP S2 results workbook.
# The new sheet is a duplicate of "Episode 10" (same headers, names and
# cell formatting) with the last column ("Felix") recomputed: every
# contestant gets an equal 1/9 share except Ricarda's row, which stays at 0
# (she already has a near-certain match elsewhere in that episode).

$wb = $excel.ActiveWorkbook
$originallyActiveSheet = $wb.ActiveSheet

$source = $wb.Worksheets.Item("Episode 10")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate the sheet (preserves styles/number formats/column layout exactly)
# and place the copy right after the last existing sheet.
$source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Episode 11"

# Update the "Felix" column (L) for the new episode. Most rows settle on an
# even 1-in-9 heat-map shade; Ricarda's row (10) keeps the neutral "no match"
# shade since she is already ~93% matched with Maurice that episode.
$oneNinth = 1.0 / 9.0
$evenColor = 0xFFE31C   # heat-map fill used for the 1/9 odds cells
$zeroColor = 0xD3D3D3   # neutral grey fill used for the 0 odds cells

$evenRows = 2,3,4,5,6,7,8,9,11
foreach ($r in $evenRows) {
    $cell = $newSheet.Range("L$r")
    $cell.Value = $oneNinth
    $cell.Interior.Color = $evenColor
}

$zeroCell = $newSheet.Range("L10")
$zeroCell.Value = 0
$zeroCell.Interior.Color = $zeroColor

# Copying a sheet makes the copy the active tab (as real Excel does); restore
# the workbook's original selection so only the sheet list actually changes.
$originallyActiveSheet.Activate()
